$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Row 2 updates
$ws.Range("D2").Value = 0.0338
$ws.Range("E2").Value = -0.033
$ws.Range("G2").Value = 0.09418977288560679
$ws.Range("H2").Value = 0.09418977288560679
$ws.Range("I2").Value = 0.07865403371882403
$ws.Range("J2").Value = 0.06500697189820717
$ws.Range("L2").Value = 0.06158562073289675
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("U2").Value = 112.3
$ws.Range("V2").Value = 0.1347977433681431
$ws.Range("W2").Value = 0.1048884670147129
$ws.Range("X2").Value = 0.0507887254447397
$ws.Range("Y2").Value = 0.05409974156997317
$ws.Range("Z2").Value = 1.740722774678632
$ws.Range("AA2").Value = 0.1131591164961031
$ws.Range("AB2").Value = 0.0487477462581859
$ws.Range("AC2").Value = 0.06441137023791715
$ws.Range("AD2").Value = 69.09999999999999
$ws.Range("AF2").Value = 69.09999999999999
$ws.Range("AG2").Value = -43.2
$ws.Range("AH2").Value = 0.07659055641764574
$ws.Range("AI2").Value = 0.06521328803322009
$ws.Range("AJ2").Value = -0.05469046714774022
$ws.Range("AK2").Value = -0.04560329357120237
$ws.Range("AL2").Value = 3.94
$ws.Range("AM2").Value = 3.94
$ws.Range("AN2").Value = 0.5215094339622641
$ws.Range("AO2").Value = 28.65482233502538
$ws.Range("AP2").Value = -0.3260377358490566
$ws.Range("AQ2").Value = 28.65482233502538
$ws.Range("D3").Value = 0.0338
$ws.Range("E3").Value = -0.033
$ws.Range("G3").Value = 0.09418977288560679
$ws.Range("H3").Value = 0.09418977288560679
$ws.Range("I3").Value = 0.07865403371882403
$ws.Range("J3").Value = 0.06500697189820717
$ws.Range("L3").Value = 0.06158562073289675
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 112.3
$ws.Range("V3").Value = 0.1347977433681431
$ws.Range("W3").Value = 0.1048884670147129
$ws.Range("X3").Value = 0.0507887254447397
$ws.Range("Y3").Value = 0.05409974156997317
$ws.Range("Z3").Value = 1.740722774678632
$ws.Range("AA3").Value = 0.1131591164961031
$ws.Range("AB3").Value = 0.0487477462581859
$ws.Range("AC3").Value = 0.06441137023791715
$ws.Range("AD3").Value = 69.09999999999999
$ws.Range("AF3").Value = 69.09999999999999
$ws.Range("AG3").Value = -43.2
$ws.Range("AH3").Value = 0.07659055641764574
$ws.Range("AI3").Value = 0.06521328803322009
$ws.Range("AJ3").Value = -0.05469046714774022
$ws.Range("AK3").Value = -0.04560329357120237
$ws.Range("AL3").Value = 3.94
$ws.Range("AM3").Value = 3.94
$ws.Range("AN3").Value = 0.5215094339622641
$ws.Range("AO3").Value = 28.65482233502538
$ws.Range("AP3").Value = -0.3260377358490566
$ws.Range("AQ3").Value = 28.65482233502538

# Clear T2:T3 (column removed from dataset)
$ws.Range("T2:T3").ClearContents()
